$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Euramet")

$ws.Range("E7").Value = 0
$ws.Range("G7").Value = -1.3734

$ws.Range("D8:H25").ClearContents()

$ws.Range("E29").Value = 0
$ws.Range("G29").Value = -1.3734

$ws.Range("E30").Value = 0
$ws.Range("F30").Value = -157.0411782090434
$ws.Range("G30").Value = -1.3734

$ws.Range("E31").Value = 0
$ws.Range("G31").Value = -1.3734

$ws.Range("E32").Value = 0
$ws.Range("F32").Value = -157.0411782090434
$ws.Range("G32").Value = -1.3734

$ws.Range("E33").Value = 0
$ws.Range("F33").Value = -157.0411782090434
$ws.Range("G33").Value = -1.3734

$ws.Range("E34").Value = 0
$ws.Range("F34").Value = -157.0411782090434
$ws.Range("G34").Value = -1.3734

$ws.Range("E35").Value = 0
$ws.Range("G35").Value = -1.3734

$ws.Range("E36").Value = 0
$ws.Range("F36").Value = -157.0411782090434
$ws.Range("G36").Value = -1.3734

$ws.Range("E37").Value = 0
$ws.Range("F37").Value = -157.0411782090434
$ws.Range("G37").Value = -1.3734

$ws.Range("E38").Value = 0
$ws.Range("F38").Value = -157.0411782090434
$ws.Range("G38").Value = -1.3734

$ws.Range("E39").Value = 0
$ws.Range("F39").Value = -157.0411782090434
$ws.Range("G39").Value = -1.3734

$ws.Range("E40").Value = 0
$ws.Range("F40").Value = -157.0411782090434
$ws.Range("G40").Value = -1.3734

$ws.Range("E41").Value = 0
$ws.Range("G41").Value = -1.3734

$ws.Range("E42").Value = 0
$ws.Range("G42").Value = -1.3734

$ws.Range("E43").Value = 0
$ws.Range("G43").Value = -1.3734

$ws.Range("E44").Value = 0
$ws.Range("F44").Value = -157.0411782090434
$ws.Range("G44").Value = -1.3734

$ws.Range("E45").Value = 0
$ws.Range("F45").Value = -157.0411782090434
$ws.Range("G45").Value = -1.3734

$ws.Range("E46").Value = 0
$ws.Range("G46").Value = -1.3734

$ws.Range("E47").Value = 0
$ws.Range("G47").Value = -1.3734
